$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.638.67'
$ws.Range("E2").Value = '  +4.13%  '

# Row 3
$ws.Range("D3").Value = '3.693.04'
$ws.Range("E3").Value = '  +6.13%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '419.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.68%  '

# Row 7
$ws.Range("D7").Value = '3.688.99'
$ws.Range("E7").Value = '  +6.34%  '

# Row 8
$ws.Range("E8").Value = '  +0.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.01%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.759'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.84%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.182'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.89%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000391'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +48.79%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.42%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.37%  '

# Row 15
$ws.Range("D15").Value = '4.283.80'
$ws.Range("E15").Value = '  +6.41%  '

# Row 16
$ws.Range("E16").Value = '  -0.55%  '

# Row 17
$ws.Range("D17").Value = '3.724.07'
$ws.Range("E17").Value = '  +7.70%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.35%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.41%  '

# Row 20
$ws.Range("E20").Value = '  +2.04%  '

# Row 21
$ws.Range("D21").Value = '66.760.53'
$ws.Range("E21").Value = '  +4.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '446.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.03%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +17.22%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.29%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.50%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.47%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.53%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.63%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.24%  '

# Row 34
$ws.Range("E34").Value = '  -1.17%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '40.81'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.62%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.23%  '

# Row 37
$ws.Range("E37").Value = '  -0.20%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0489'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.07%  '

# Row 39
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +39.61%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0739'
$ws.Range("E40").Value = '  +11.54%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.149'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.18%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '29.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +32.84%  '

# Row 43
$ws.Range("E43").Value = '  +0.12%  '

# Row 44
$ws.Range("E44").Value = '  +1.26%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.75%  '

# Row 46
$ws.Range("E46").Value = '  +3.64%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.17%  '

# Row 48
$ws.Range("E48").Value = '  -3.19%  '

# Row 49
$ws.Range("B49").Value = 'TheGraph'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.305'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.83%  '

# Row 50
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.16%  '

# Row 51
$ws.Range("E51").Value = '  +16.39%  '
